$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (Sending, Ligand, Receptor, Target,
# then 16 numeric columns E..T), reflecting the corrected NATMI run that
# adds the sCs cluster into the 3x3 ECs/FAPs/sCs grid for Bmp4-Bmpr1b.
$data = @(
    @(2, "ECs", "Bmp4", "Bmpr1b", "ECs", 3, 1, 8.675694999999999, 26.027085, 0.5592117158070719, 0.5592117158070719, 1, 0.3333333333333333, 0.039133, 0.117399, 0.00872147833139909, 0.00872147833139909, 0.339505972435, 3.055553751915, 0.004877152862075883, 0.004877152862075883),
    @(3, "ECs", "Bmp4", "Bmpr1b", "FAPs", 3, 1, 8.675694999999999, 26.027085, 0.5592117158070719, 0.5592117158070719, 3, 1, 4.187332, 12.561996, 0.933220691088698, 0.9332206910886979, 36.32801529574, 326.95213766166, 0.5218679438903722, 0.5218679438903721),
    @(4, "ECs", "Bmp4", "Bmpr1b", "sCs", 3, 1, 8.675694999999999, 26.027085, 0.5592117158070719, 0.5592117158070719, 3, 1, 0.2605036666666667, 0.7815110000000001, 0.05805783057990302, 0.05805783057990302, 2.260050358381667, 20.340453225435, 0.03246661905462386, 0.03246661905462386),
    @(5, "FAPs", "Bmp4", "Bmpr1b", "ECs", 3, 1, 6.316050666666666, 18.948152, 0.407115456505913, 0.407115456505913, 1, 0.3333333333333333, 0.039133, 0.117399, 0.00872147833139909, 0.00872147833139909, 0.2471660107386667, 2.224494096648, 0.003550648632293968, 0.003550648632293968),
    @(6, "FAPs", "Bmp4", "Bmpr1b", "FAPs", 3, 1, 6.316050666666666, 18.948152, 0.407115456505913, 0.407115456505913, 3, 1, 4.187332, 12.561996, 0.933220691088698, 0.9332206910886979, 26.44740107015467, 238.026609631392, 0.3799285676733389, 0.3799285676733388),
    @(7, "FAPs", "Bmp4", "Bmpr1b", "sCs", 3, 1, 6.316050666666666, 18.948152, 0.407115456505913, 0.407115456505913, 3, 1, 0.2605036666666667, 0.7815110000000001, 0.05805783057990302, 0.05805783057990302, 1.645354357519111, 14.808189217672, 0.02363624020028017, 0.02363624020028017),
    @(8, "sCs", "Bmp4", "Bmpr1b", "ECs", 2, 0.6666666666666666, 0.5224053333333334, 1.567216, 0.03367282768701513, 0.03367282768701512, 1, 0.3333333333333333, 0.039133, 0.117399, 0.00872147833139909, 0.00872147833139909, 0.02044328790933334, 0.183989591184, 0.0002936768370292378, 0.0002936768370292377),
    @(9, "sCs", "Bmp4", "Bmpr1b", "FAPs", 2, 0.6666666666666666, 0.5224053333333334, 1.567216, 0.03367282768701513, 0.03367282768701512, 3, 1, 4.187332, 12.561996, 0.933220691088698, 0.9332206910886979, 2.187484569237334, 19.687361123136, 0.0314241795249869, 0.03142417952498689),
    @(10, "sCs", "Bmp4", "Bmpr1b", "sCs", 2, 0.6666666666666666, 0.5224053333333334, 1.567216, 0.03367282768701513, 0.03367282768701512, 3, 1, 0.2605036666666667, 0.7815110000000001, 0.05805783057990302, 0.05805783057990302, 0.1360885048195556, 1.224796543376, 0.001954971324998992, 0.001954971324998991)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    for ($i = 5; $i -le 20; $i++) {
        $ws.Cells.Item($r, $i).Value = $row[$i]
    }
}

Write-Host "done"